$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-12-10"

# Update the "2022 (through 12-09)" column header label
$ws.Range("I1").Value = "2022 (through 12-10)"

# Update December 2022 count and the recalculated Total
$ws.Range("I13").Value = 44
$ws.Range("I14").Value = 1560
